$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow edits, then re-protect afterward.
$ws.Unprotect()

# Update the confidential/date disclaimer text (A18): 2021-05-20 -> 2021-05-21
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15

$ws.Range("D2").Value = 0.05660989403245811
$ws.Range("E2").Value = -0.0005579059928402241

$ws.Range("D3").Value = 0.02388595905395758
$ws.Range("E3").Value = -0.0001928268414963519

$ws.Range("D4").Value = 0.03116649833366781
$ws.Range("E4").Value = -0.01064031920957631

$ws.Range("D5").Value = 0.03255417635815121
$ws.Range("E5").Value = 0.002304147465437945

$ws.Range("D6").Value = 0.0370863367339889
$ws.Range("E6").Value = -0.007025761124121677

$ws.Range("D7").Value = 0.01870437325841144
$ws.Range("E7").Value = 0.006156119182467279

$ws.Range("D8").Value = 0.00444728106566303
$ws.Range("E8").Value = 0.002219263204615851

$ws.Range("D9").Value = 0.00682850760459401
$ws.Range("E9").Value = 0.006937752938909192

$ws.Range("D10").Value = 0.07370682423074014
$ws.Range("E10").Value = 0.0005356186395286677

$ws.Range("D11").Value = 0.07382526047749548
$ws.Range("E11").Value = 0.0005347593582889498

$ws.Range("D12").Value = 0.1444711657087587
$ws.Range("E12").Value = 0.003206295999417019

$ws.Range("D13").Value = 0.3831195449416208
$ws.Range("E13").Value = 0.0003503547341683344

$ws.Range("D14").Value = 0.1135941782004928
$ws.Range("E14").Value = 0.004622308721566881

$ws.Range("E15").Value = 0.0008204999983552064

# Re-protect the sheet (password unknown/irrecoverable from the legacy hash;
# re-apply protection so the sheet remains protected as in the source file).
$ws.Protect()
